# New columns I ("I0") and J ("IF") added, mirroring the existing stat
# columns (B..H) both in data and in header formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the two new header cells the same style as the existing header
# cells (bold, centered, bordered) by copying H1's formatting onto them.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Per-row data for the new I0 / IF columns.
$iValues = @(1, 1, 8, 4, 1, 1, 1, 1, 1, 4, 1, 1, 1)
$jValues = @(5, 5, 9, 6, 4, 2, 5, 5, 3, 6, 3, 3, 2)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
